$p = $ppt.ActivePresentation

# --- Slide 2: "Executive Program Status" content placeholder ---
$s2 = $p.Slides.Item(2)
$shape2 = $s2.Shapes.Item(2)

$slide2Paragraphs = @(
    '🔴',
    '## Overall Health & Status',
    '*   **Persistent Under-delivery:** The team consistently delivers only 50-60% of committed work, with recent sprints frequently falling below 40%, severely impacting predictability and roadmap execution.',
    '*   **Critical Quality Decline:** Defect density is alarmingly high, with recent sprints averaging over one defect per story, indicating significant quality issues and potential for extensive rework.',
    '*   **Overloaded Team & Bottlenecks:** Key team members are significantly over-capacity, while a large volume of work is stalled in SIT and UAT environments, indicating severe pipeline blockages and an unsustainable workload.',
    '*   **Unmanaged Critical Risks:** A substantial number of high-impact risks, issues, and dependencies remain open and unmitigated, posing significant threats to project timelines and overall stability.',
    '',
    '## Key Risks and Impediments',
    '*   **High External Dependency Risk:** A vast majority (12 of 14) of identified dependencies are open and unmanaged, representing a critical external blocking factor severely impacting delivery timelines.',
    '*   **Team Burnout & Attrition:** Consistent over-capacity for multiple team members combined with the pressure of high defect rates creates a high risk of burnout, reduced morale, and potential key talent attrition.',
    '*   **Compromised Product Quality & Rework:** The persistent high defect density through SIT/UAT stages significantly increases the likelihood of critical bugs reaching production, escalating support costs and damaging customer trust.',
    '*   **Unpredictable Delivery & Scope Creep:** The lack of consistent velocity and high number of open issues contribute to a highly unpredictable delivery cadence, making accurate forecasting and effective scope management extremely challenging.',
    '',
    '## Executive Recommendations',
    '*   **Immediate RAID Prioritization & Mitigation:** Establish an immediate, cross-functional task force to review and actively mitigate all critical open risks, issues, and especially dependencies, with daily tracking and reporting.',
    '*   **Re-evaluate Capacity & Workload Balancing:** Conduct an urgent review of team capacity, redistribute workload more equitably, and consider pausing new commitments to stabilize the team and allow for resolution of technical debt.',
    '*   **Enhance Quality Gates & Definition of Done:** Implement stricter quality gates upstream (e.g., improved unit testing, peer reviews, shift-left testing) and re-enforce the Definition of Done to reduce defect escape rates into later stages.',
    '*   **Right-size Sprint Commitments:** Adjust future sprint commitments to reflect actual historical velocity (e.g., targeting 50-60% of current commitments) to build predictability and enable the team to focus on quality and resolution of current impediments.'
)

$shape2.TextFrame.TextRange.Text = [string]::Join("`r", $slide2Paragraphs)

# --- Slide 3: "Velocity & Completion Trend" content placeholder ---
$s3 = $p.Slides.Item(3)
$shape3 = $s3.Shapes.Item(2)

$slide3FirstParagraphs = @(
    '❌',
    '*   Significant gap exists between committed and completed points.',
    '*   Average completion rate is consistently low, around 55%.',
    '*   Completed velocity shows high variability and frequent dips.',
    '*   Team struggles with predictability and realistic sprint commitments.'
)

$slide3LastParagraphLines = @(
    '❌',
    '*   Overall low sprint goal completion.',
    '*   Average completion rate is approximately 54%.',
    '*   Completion rates fluctuate widely (33% to 74%).',
    '*   No sprint achieved full goal completion.',
    '*   Suggests issues with planning, estimation, or delivery.'
)

$slide3Text = [string]::Join("`r", $slide3FirstParagraphs) + "`r" + [string]::Join("`v", $slide3LastParagraphLines)

$shape3.TextFrame.TextRange.Text = $slide3Text
